$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add "?" marker (H column) to a few existing rows that previously lacked it
$ws.Range("H2").Value = "?"
$ws.Range("H4").Value = "?"
$ws.Range("H5").Value = "?"
$ws.Range("H6").Value = "?"

# Clear the father/mother-link values that had been on row 7 (J7/K7) --
# they belong to the new person being inserted as row 8 instead.
$ws.Range("J7").ClearContents()
$ws.Range("K7").ClearContents()

# Insert a new row at 8 (shifts old rows 8-36 down to 9-37), copying
# formatting from the row above as Excel normally does.
$ws.Rows(8).Insert()

# Populate the newly inserted row 8 with the new person's data.
$ws.Range("A8").Value = 140003
$ws.Range("B8").Value = 14
$ws.Range("C8").Value = "x"
$ws.Range("D8").Value = "Dương Danh Mô"
$ws.Range("E8").Value = "?"
$ws.Range("F8").Value = "?"
$ws.Range("H8").Value = "?"
$ws.Range("I8").Value = "?"

# Update the selection to match the author's saved cursor position.
$ws.Range("J7:K8").Select() | Out-Null
